$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Dyson Daniels'
$ws.Cells.Item(2, 2).Value = 'PG,SG,SF'
$ws.Cells.Item(2, 3).Value = 'Atlanta Hawks'
$ws.Cells.Item(3, 1).Value = 'Malik Beasley'
$ws.Cells.Item(3, 2).Value = 'SG,SF'
$ws.Cells.Item(3, 3).Value = 'Detroit Pistons'
$ws.Cells.Item(4, 1).Value = 'Cam Thomas'
$ws.Cells.Item(4, 2).Value = 'SG,SF'
$ws.Cells.Item(4, 3).Value = 'Brooklyn Nets'
$ws.Cells.Item(5, 1).Value = 'Toumani Camara'
$ws.Cells.Item(5, 2).Value = 'SG,SF,PF'
$ws.Cells.Item(5, 3).Value = 'Portland Trail Blazers'
$ws.Cells.Item(6, 1).Value = 'Jaden McDaniels'
$ws.Cells.Item(6, 2).Value = 'SF,PF'
$ws.Cells.Item(6, 3).Value = 'Minnesota Timberwolves'
$ws.Cells.Item(7, 1).Value = 'Julius Randle'
$ws.Cells.Item(7, 2).Value = 'PF,C'
$ws.Cells.Item(7, 3).Value = 'Minnesota Timberwolves'
$ws.Cells.Item(8, 1).Value = 'De''Andre Hunter'
$ws.Cells.Item(8, 2).Value = 'SF,PF'
$ws.Cells.Item(8, 3).Value = 'Cleveland Cavaliers'
$ws.Cells.Item(9, 1).Value = 'Michael Porter Jr.'
$ws.Cells.Item(9, 2).Value = 'SF,PF'
$ws.Cells.Item(9, 3).Value = 'Denver Nuggets'
$ws.Cells.Item(10, 1).Value = 'Alperen Sengün'
$ws.Cells.Item(10, 2).Value = 'C'
$ws.Cells.Item(10, 3).Value = 'Houston Rockets'
$ws.Cells.Item(11, 1).Value = 'Guerschon Yabusele'
$ws.Cells.Item(11, 2).Value = 'PF,C'
$ws.Cells.Item(11, 3).Value = 'Philadelphia 76ers'
$ws.Cells.Item(12, 1).Value = 'Josh Hart'
$ws.Cells.Item(12, 2).Value = 'SG,SF,PF'
$ws.Cells.Item(12, 3).Value = 'New York Knicks'
$ws.Cells.Item(13, 1).Value = 'Donovan Mitchell'
$ws.Cells.Item(13, 2).Value = 'PG,SG'
$ws.Cells.Item(13, 3).Value = 'Cleveland Cavaliers'
$ws.Cells.Item(14, 1).Value = 'Devin Vassell'
$ws.Cells.Item(14, 2).Value = 'SG,SF'
$ws.Cells.Item(14, 3).Value = 'San Antonio Spurs'
$ws.Cells.Item(15, 1).Value = 'Kelly Oubre Jr.'
$ws.Cells.Item(15, 2).Value = 'SG,SF'
$ws.Cells.Item(15, 3).Value = 'Philadelphia 76ers'
$ws.Cells.Item(16, 1).Value = 'Kelly Olynyk'
$ws.Cells.Item(16, 2).Value = 'PF,C'
$ws.Cells.Item(16, 3).Value = 'New Orleans Pelicans'
$ws.Cells.Item(17, 1).Value = 'Domantas Sabonis'
$ws.Cells.Item(17, 2).Value = 'C'
$ws.Cells.Item(17, 3).Value = 'Sacramento Kings'
$ws.Cells.Item(18, 1).Value = 'Kristaps Porzingis'
$ws.Cells.Item(18, 2).Value = 'PF,C'
$ws.Cells.Item(18, 3).Value = 'Boston Celtics'
